# Update "alquiler_ajustado_ponderaciones_IPC_2018" sheet:
#   - C44 (the 041 Alquiler de vivienda adjusted weight) changes from 250 -> 226.58
#   - C2:C43 are re-derived so that column C (rows 2-43) is column B (rows 2-43)
#     rescaled to sum to (1000 - C44), keeping the total of C2:C44 at 1000.
#   - B27 is re-written with its (floating point re-evaluated) value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow  = 43
$rentRow  = 44

# New adjusted weight for "041 Alquiler de vivienda" (row 44, column C).
$newRentWeight = 226.58
$ws.Cells.Item($rentRow, 3).Value = $newRentWeight

# B27 is re-entered (tiny floating point re-evaluation noise vs. the original 0.72).
$ws.Cells.Item(27, 2).Value = 0.7199999999999999

# Sum column B across rows 2..43 using the (possibly just-updated) cell values.
$sumB = 0
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $sumB = $sumB + $ws.Cells.Item($r, 2).Value()
}

# Remaining weight to distribute across rows 2..43 so the full column (2..44) sums to 1000.
$targetSum = 1000 - $newRentWeight

# Re-derive each C value proportionally to its B weight.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $b = $ws.Cells.Item($r, 2).Value()
    $ws.Cells.Item($r, 3).Value = $b * $targetSum / $sumB
}
